$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: assign a text value to a cell while preserving its original
# style/number-format identity, even when the text looks like a number
# (Excel would otherwise silently convert it to a numeric cell).
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '28.057.18'
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').Value = '1.890.15'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue $ws.Range('D5') '313.82'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('E6').Value = '  +0.06%  '
Set-TextValue $ws.Range('D7') '0.5046'
$ws.Range('E7').Value = '  -1.88%  '
Set-TextValue $ws.Range('D8') '0.3893'
$ws.Range('E8').Value = '  -1.93%  '
Set-TextValue $ws.Range('D9') '0.09243'
$ws.Range('E9').Value = '  -5.85%  '
Set-TextValue $ws.Range('D10') '1.126'
$ws.Range('E10').Value = '  -3.03%  '
Set-TextValue $ws.Range('D11') '41.83'
$ws.Range('E11').Value = '  -1.11%  '
Set-TextValue $ws.Range('D12') '6.385'
$ws.Range('E12').Value = '  -2.54%  '
Set-TextValue $ws.Range('D13') '20.79'
$ws.Range('E13').Value = '  -2.46%  '
$ws.Range('D14').Value = '1.899.11'
$ws.Range('E14').Value = '  -0.74%  '
Set-TextValue $ws.Range('D15') '7.283'
$ws.Range('E15').Value = '  -4.17%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('E18').Value = '  -3.38%  '
Set-TextValue $ws.Range('D19') '0.06655'
$ws.Range('E19').Value = '  -0.11%  '
Set-TextValue $ws.Range('D20') '17.82'
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('E21').Value = '  +0.09%  '
Set-TextValue $ws.Range('D22') '6.208'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').Value = '28.109.96'
$ws.Range('E23').Value = '  -1.88%  '
Set-TextValue $ws.Range('D24') '11.39'
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('D26').Value = '2.114.92'
$ws.Range('E26').Value = '  -0.79%  '
Set-TextValue $ws.Range('D27') '2.540'
$ws.Range('E27').Value = '  -6.21%  '
Set-TextValue $ws.Range('D28') '158.42'
$ws.Range('E28').Value = '  -0.96%  '
Set-TextValue $ws.Range('D29') '20.80'
$ws.Range('E29').Value = '  -2.33%  '
Set-TextValue $ws.Range('D30') '126.77'
Set-TextValue $ws.Range('D31') '1.077'
$ws.Range('E31').Value = '  -2.73%  '
Set-TextValue $ws.Range('D32') '0.1055'
$ws.Range('E32').Value = '  -2.59%  '
Set-TextValue $ws.Range('D33') '5.605'
$ws.Range('E33').Value = '  -2.84%  '
Set-TextValue $ws.Range('D34') '3.609'
$ws.Range('E34').Value = '  -0.89%  '
Set-TextValue $ws.Range('D35') '9.465'
$ws.Range('E35').Value = '  -4.13%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D36') '1.348'
$ws.Range('E36').Value = '  +12.80%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D37') '0.06590'
$ws.Range('E37').Value = '  -3.41%  '
Set-TextValue $ws.Range('D38') '0.02404'
$ws.Range('E38').Value = '  -1.54%  '
Set-TextValue $ws.Range('D39') '0.2195'
$ws.Range('E39').Value = '  -1.90%  '
Set-TextValue $ws.Range('D40') '1.217'
$ws.Range('E40').Value = '  -4.27%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D41') '11.66'
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D42') '0.6430'
$ws.Range('E42').Value = '  -0.37%  '
Set-TextValue $ws.Range('D43') '4.962'
$ws.Range('E43').Value = '  -3.34%  '
$ws.Range('E44').Value = '  +0.06%  '
Set-TextValue $ws.Range('D45') '13.28'
$ws.Range('E45').Value = '  -2.65%  '
Set-TextValue $ws.Range('D46') '0.6039'
$ws.Range('E46').Value = '  -1.02%  '
Set-TextValue $ws.Range('D47') '1.300'
$ws.Range('E47').Value = '  +1.63%  '
Set-TextValue $ws.Range('D48') '3.689'
$ws.Range('E48').Value = '  -2.58%  '
Set-TextValue $ws.Range('D49') '1.999'
$ws.Range('E49').Value = '  -2.08%  '
Set-TextValue $ws.Range('D50') '122.15'
$ws.Range('E50').Value = '  -2.68%  '
$ws.Range('E51').Value = '  -1.97%  '
